$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target data for rows 2-55 (journal name variants -> normalized name), per the commit diff
$data = @(
    @("ACM TRANSACTIONS ON GRAPHICS", "ACM Transactions on Graphics"),
    @("ANNALS OF OPERATION RESEARCH", "RAIRO-OPERATIONS RESEARCH"),
    @("COMPUTER NETWORKS (1999)", "Computer Networks (1999)"),
    @("IEEE COMMUNICATIONS MAGAZINE (PRINT)", "IEEE Communications Magazine (Print)"),
    @("IEEE COMMUNICATIONS SURVEYS AND TUTORIALS", "IEEE Communications Surveys and Tutorials"),
    @("IEEE NETWORK", "IEEE Network"),
    @("IEEE TRANSACTIONS ON COMMUNICATIONS (PRINT)", "IEEE Transactions on Cognitive Communications and Networking"),
    @("IEEE TRANSACTIONS ON NUCLEAR SCIENCE", "IEEE Transactions on Nanobioscience"),
    @("IEEE WIRELESS COMMUNICATIONS", "IEEE Wireless Communications"),
    @("IEEE/ACM TRANSACTIONS ON NETWORKING (PRINT)", "IEEE/ACM Transactions on Networking"),
    @("INTEGRATED COMPUTER-AIDED ENGINEERING", "Integrated Computer-Aided Engineering"),
    @("SENSORS (BASEL)", "SENSORS"),
    @("SIGMOD RECORD", "SIGMOD Record"),
    @("JOURNAL OF PARALLEL AND DISTRIBUTED COMPUTING (PRINT)", "Journal of Parallel and Distributed Computing (Print)"),
    @("OPTIMIZATION LETTERS (PRINT)", "Optimization Letters (Print)"),
    @("PARALLEL COMPUTING", "Parallel Computing"),
    @("PERFORMANCE EVALUATION", "Performance Evaluation"),
    @("SCIENCE OF COMPUTER PROGRAMMING (PRINT)", "SCIENCE OF COMPUTER PROGRAMMING"),
    @("ANNALES DES TELECOMMUNICATIONS", "Annals of Telecommunications"),
    @("AUTOMATED SOFTWARE ENGINEERING", "Automated Software Engineering"),
    @("COMPUTERS IN ENTERTAINMENT (CIE)", "Computers in Entertainment : CIE"),
    @("EARTH SCIENCE INFORMATICS (PRINT)", "Earth Science Informatics"),
    @("INTERNATIONAL JOURNAL OF BUSINESS PROCESS INTEGRATION AND MANAGEMENT (PRINT)", "International Journal of Business Process Integration and Management (Print)"),
    @("INTERNATIONAL JOURNAL OF COMMUNICATION SYSTEMS (PRINT)", "International Journal of Computational Fluid Dynamics (Print)"),
    @("INTERNATIONAL JOURNAL OF DISTRIBUTED SENSOR NETWORKS", "International Journal of Distributed Sensor Networks"),
    @("INTERNATIONAL JOURNAL OF SOFTWARE ENGINEERING AND KNOWLEDGE ENGINEERING", "International Journal of Software Engineering and Knowledge Engineering"),
    @("JOURNAL OF COMPUTATIONAL SCIENCE", "Journal of Computational Science"),
    @("JOURNAL OF INTERNET SERVICES AND APPLICATIONS", "Journal of Internet Services and Applications"),
    @("JOURNAL OF UNIVERSAL COMPUTER SCIENCE (PRINT)", "Journal of Universal Computer Science (Print)"),
    @("NATURAL COMPUTING", "Astronomy and Computing"),
    @("SOFTWARE QUALITY JOURNAL", "ProQuality (UFLA)"),
    @("INTERNATIONAL JOURNAL OF INNOVATIVE COMPUTING AND APPLICATIONS (ONLINE)", "International Journal of Innovative Computing and Applications (Online)"),
    @("APPLICABLE ANALYSIS AND DISCRETE MATHEMATICS", "Electronic Notes in Discrete Mathematics"),
    @("IEEE PHOTONICS JOURNAL", "IEEE Potentials"),
    @("JOURNAL OF MATHEMATICAL ANALYSIS AND APPLICATIONS (PRINT)", "JOURNAL OF MATHEMATICAL MODELLING AND ALGORITHMS (PRINT)"),
    @("REVISTA DE INFORMÁTICA TEÓRICA E APLICADA (IMPRESSO)", "Revista de Informática Teórica e Aplicada (Impresso)"),
    @("COMUNICAR (HUELVA)", "INFOCOMP (UFLA)"),
    @("AÇÃO ERGONÔMICA", "Espaço Energia"),
    @("CLEI ELECTRONIC JOURNAL", "CLEI Electronic Journal"),
    @("JOURNAL OF COMMUNICATION AND INFORMATION SYSTEMS (ONLINE)", "Journal of Communication and Information Systems (Online)"),
    @("JOURNAL OF COMPUTACIONAL INTERDISCIPLINARY SCIENCES", "Journal of Computational Interdisciplinary Sciences"),
    @("RESI. REVISTA ELETRÔNICA DE SISTEMAS DE INFORMAÇÃO", "RESI. Revista Eletrônica de Sistemas de Informação"),
    @("REVISTA BRASILEIRA DE REDES DE COMPUTADORES E SISTEMAS DISTRIBUÍDOS", "Revista Brasileira de Redes de Computadores e Sistemas Distribuídos"),
    @("SOCIAL NETWORK ANALYSIS AND MINING", "Statistical Analysis and Data Mining (Online)"),
    @("AACE INTERNATIONAL TRANSACTIONS", "IEEE Latin America Transactions"),
    @("EM DEBATE (BELO HORIZONTE)", "Fonte (Belo Horizonte)"),
    @("INTERNATIONAL JOURNAL OF NETWORKS AND COMMUNICATIONS", "Cluster Computing-The Journal of Networks Software Tools and Applications"),
    @("INTERNATIONAL JOURNAL OF SCIENTIFIC RESEARCH AND REVIEWS", "International Journal of Open Source Software and Processes"),
    @("LECTURE NOTES IN COMPUTER SCIENCE", "Lecture Notes in Computer Science"),
    @("PROCEDIA COMPUTER SCIENCE", "Procedia Computer Science"),
    @("REVISTA ELETRÔNICA DE INICIAÇÃO CIENTÍFICA", "REIC. Revista Eletrônica de Iniciação Científica"),
    @("REVISTA TECNOLOGIAS NA EDUCAÇÃO", "Revista Tecnologia da Informação"),
    @("S & G. SISTEMAS & GESTÃO", "S & G. Sistemas & Gestão"),
    @("TECNOLOGIA EDUCACIONAL", "Tecnologia Educacional")
)

$startRow = 2
$endRow = $startRow + $data.Length - 1

# Extend formatting (bold, centered, bordered style used by existing column-A cells)
# from the last currently-formatted row (A37) down through the new rows (A38:A55),
# so the new rows match the look of the existing ones.
$ws.Range("A2").Copy($ws.Range("A38:A55")) | Out-Null

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $pair = $data[$i]
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
}

Write-Host ("Wrote rows {0}-{1}" -f $startRow, $endRow)
